$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) keeps its text formatting so numeric-looking
# strings like "0.9996" or "1.0000" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "22.069.14"
$ws.Range("E2").Value = "  -1.46%  "
$ws.Range("D3").Value = "1.557.30"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "0.9996"
$ws.Range("D6").Value = "287.63"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "0.3869"
$ws.Range("E7").Value = "  +4.51%  "
$ws.Range("D8").Value = "0.3243"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").Value = "43.11"
$ws.Range("E9").Value = "  -7.29%  "
$ws.Range("D10").Value = "1.123"
$ws.Range("E10").Value = "  -1.85%  "
$ws.Range("D11").Value = "0.07367"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "1.0000"
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").Value = "19.36"
$ws.Range("E13").Value = "  -5.19%  "
$ws.Range("D14").Value = "5.701"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "6.808"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.00001126"
$ws.Range("E16").Value = "  +2.50%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.555.90"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "0.06609"
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("D19").Value = "85.25"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "6.390"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "0.9990"
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "15.99"
$ws.Range("E22").Value = "  -1.50%  "
$ws.Range("D23").Value = "11.48"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").Value = "22.076.77"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").Value = "2.320"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").Value = "2.560"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("D27").Value = "149.51"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").Value = "18.87"
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("D29").Value = "4.868"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "1.732.72"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "120.85"
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("D32").Value = "1.112"
$ws.Range("E32").Value = "  +6.60%  "
$ws.Range("D33").Value = "5.855"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "1.704"
$ws.Range("E34").Value = "  -13.41%  "
$ws.Range("D35").Value = "0.08195"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "9.282"
$ws.Range("E36").Value = "  -3.90%  "
$ws.Range("D37").Value = "0.06241"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("D38").Value = "0.02299"
$ws.Range("E38").Value = "  -3.53%  "
$ws.Range("D39").Value = "5.232"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  -3.48%  "
$ws.Range("D41").Value = "1.221"
$ws.Range("E41").Value = "  -5.98%  "
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").Value = "0.9990"
$ws.Range("D44").Value = "0.5953"
$ws.Range("E44").Value = "  -2.41%  "
$ws.Range("D45").Value = "13.58"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "3.717"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("D47").Value = "0.5754"
$ws.Range("E47").Value = "  -3.13%  "
$ws.Range("D48").Value = "1.928"
$ws.Range("E48").Value = "  -3.87%  "
$ws.Range("D49").Value = "119.08"
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("D50").Value = "1.159"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "0.06892"
$ws.Range("E51").Value = "  -3.65%  "
